$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

$ws.Range("A56").Value = "The Great Influenza"
$ws.Range("B56").Value = "John Barry;Scott Brick"

$ws.Range("C2").Copy($ws.Range("C56"))
$ws.Range("C56").Value = 43935

$ws.Range("D2").Copy($ws.Range("D56"))
$ws.Range("D56").Value = 43939

$ws.Range("E56").Value = "influenza;pandemic;science;disease;history;spanish flu"
$ws.Range("F56").Value = "Audio"
$ws.Range("G56").Value = "19 Hours 34 Mins"

$ws.Range("A57").Select()
